$wb = $excel.ActiveWorkbook

# Add the new "ODI Batting Extra" worksheet after the last existing sheet
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "ODI Batting Extra"

# Header row (bold, thin border all sides, centered horizontal/top vertical -
# matches the "header" style used on the other sheets)
$header = $ws.Range("A1:F1")
$header.Font.Bold = $true
$header.Borders.LineStyle = 1
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4160

$ws.Range("A1").Value = "MATCH_CODE"
$ws.Range("B1").Value = "BATTING_POSITION"
$ws.Range("C1").Value = "NUM_4"
$ws.Range("D1").Value = "NUM_6"
$ws.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$ws.Range("F1").Value = "MAN_OF_MATCH"

# Data row
# MATCH_CODE is stored as text (not a number)
$ws.Range("A2").Value = "'4641"
$ws.Range("A2").Style = "Normal"

# BATTING_POSITION is a real number
$ws.Range("B2").Value = 1

# NUM_4 / NUM_6 are text "0"
$ws.Range("C2").Value = "'0"
$ws.Range("C2").Style = "Normal"
$ws.Range("D2").Value = "'0"
$ws.Range("D2").Style = "Normal"

# PERCENT_RUNS_OF_TOTAL is the literal text "0.97%"
$ws.Range("E2").Value = "'0.97%"
$ws.Range("E2").Style = "Normal"

# MAN_OF_MATCH text
$ws.Range("F2").Value = "NO"

# Keep the originally active tab (Player Info) selected, like before the edit
$wb.Worksheets.Item(1).Activate()
